# #31 Added InTableOptionDefinition. Added readContext methods in XlBeanReader.
# Adds a new "inTableOptions" worksheet (after "offset") exercising the
# in-table option syntax, mirroring the existing "limit"/"offset" sheets.

$wb = $excel.ActiveWorkbook

# --- add the new sheet at the end of the workbook, after "offset" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "inTableOptions"

# --- header row ---
$ws.Range("A1").Value = "####"
$ws.Range("C1").Value = "single"
$ws.Range("E1").Value = "optionsInTable#field1"
$ws.Range("F1").Value = "optionsInTable#field2"

# --- "single" option block (rows 3-6) ---
$ws.Range("A3").Value = "single?type"
$ws.Range("C3").Value = "string"

$ws.Range("A4").Value = "single?custom"
$ws.Range("C4").Value = "test value for custom option"

$ws.Range("A6").Value = "single"
$ws.Range("C6").Value = 100

# --- "optionsInTable" block (rows 8-9) ---
$ws.Range("A8").Value = "optionsInTable?type"
$ws.Range("E8").Value = "string"

$ws.Range("A9").Value = "optionsInTable?customType"
$ws.Range("E9").Value = "hoge"
$ws.Range("F9").Value = "fuga"

# --- "optionsInTable" data table (rows 11-18) ---
$ws.Range("A11").Value = "optionsInTable#~"

for ($r = 11; $r -le 18; $r++) {
    $ws.Range("E$r").Value = $r - 10
    $ws.Range("F$r").Formula = "=E$r+100"
}

# --- blank-but-styled cells (shaded background carries through even though empty) ---
$ws.Range("C2").Value = ""
$ws.Range("C5").Value = ""
$ws.Range("B3").Value = ""
$ws.Range("B4").Value = ""
$ws.Range("B6").Value = ""
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = ""
$ws.Range("F8").Value = ""
$ws.Range("B9").Value = ""
$ws.Range("C9").Value = ""
$ws.Range("D9").Value = ""
$ws.Range("B11").Value = ""
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = ""
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = ""
$ws.Range("E10").Value = ""
$ws.Range("F10").Value = ""

# --- shading: column C (rows 1-6) + the A:B "single" option labels, light orange (theme Accent4, lighter 80%) ---
$ws.Range("C1:C6").Interior.ThemeColor = 8
$ws.Range("A3:B4").Interior.ThemeColor = 8
$ws.Range("A6:B6").Interior.ThemeColor = 8
# --- shading: columns E:F (rows 1-18) and A8:D9 / A11:D11 light gold (theme Accent2, lighter 80%) ---
$ws.Range("E1:F18").Interior.ThemeColor = 6
$ws.Range("A8:D9").Interior.ThemeColor = 6
$ws.Range("A11:D11").Interior.ThemeColor = 6

# --- column widths (best-fit) ---
$ws.Columns.Item(1).ColumnWidth = 25.7
$ws.Range("E1:F1").ColumnWidth = 21.7

# --- selection / active cell matches the authored file ---
$ws.Range("C5").Select()

$wb.Worksheets.Item("offset").Select()
$ws.Activate()
